# Auto update: 2025-12-05 02:39:10
# Refresh the daily quantum-stock scoring table on Sheet1.
#
# Row order changes: the "International Business Machines" / "IBM" entry
# moves up from row 5 to row 4, and "Rigetti Computing, Inc." / "RGTI"
# moves down from row 4 to row 5. All numeric metrics are refreshed with
# the latest run's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order the stock name / ticker for rows 4 and 5 ---
$ws.Range("B4").Value = "International Business Machines"
$ws.Range("C4").Value = "IBM"
$ws.Range("B5").Value = "Rigetti Computing, Inc."
$ws.Range("C5").Value = "RGTI"

# --- Row 2 (IonQ, Inc. / IONQ): 종가, RSI, 5일수익률, 최종점수, MACRO_SCORE ---
$ws.Range("D2").Value = 54.46
$ws.Range("E2").Value = 64.90000000000001
$ws.Range("F2").Value = 16.11
$ws.Range("K2").Value = 58.1
$ws.Range("N2").Value = 53.62998959737769

# --- Row 3 (D-Wave Quantum Inc. / QBTS) ---
$ws.Range("D3").Value = 28.21
$ws.Range("E3").Value = 64.8
$ws.Range("F3").Value = 25.9
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 57.5
$ws.Range("N3").Value = 53.62998959737769

# --- Row 4 (now International Business Machines / IBM) ---
$ws.Range("D4").Value = 306.89
$ws.Range("E4").Value = 51.9
$ws.Range("F4").Value = 1.21
$ws.Range("H4").Value = 66
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 56.3
$ws.Range("N4").Value = 53.62998959737769

# --- Row 5 (now Rigetti Computing, Inc. / RGTI) ---
$ws.Range("D5").Value = 29.01
$ws.Range("E5").Value = 60.9
$ws.Range("F5").Value = 13.45
$ws.Range("H5").Value = 63
$ws.Range("J5").Value = 83
$ws.Range("K5").Value = 55.1
$ws.Range("N5").Value = 53.62998959737769
